$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new data row (Agriculture & Co-Operation dept / Annadata Sukhibhava) ---

# A4: Source_Name (reuse "Andhra Pradesh Finance Department"), bordered like A2
$ws.Range("A4").Value = "Andhra Pradesh Finance Department"
$ws.Range("A4").Borders.ColorIndex = 1
$ws.Range("A4").Borders.LineStyle = 1

# B4: Source_URL + hyperlink
$ws.Range("B4").Value = "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-11.pdf"
$ws.Range("B4").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("B4"), "https://apfinance.gov.in/...Bud@et24-25/documents/Volume-III-11.pdf") | Out-Null

# C4: Document_Name
$ws.Range("C4").Value = "Agriculture & Co-Operation & Food, Civil Supplies & Consumers Affairs Department"
$ws.Range("C4").WrapText = $true

# D4: Financial Year(s) covered (reuse "2024-25")
$ws.Range("D4").Value = "2024-25"

# E4: Date Downloaded
$ws.Range("E4").Value = "7/7/2025"
$ws.Range("E4").NumberFormat = "mm-dd-yy"

# F4: Method of Acquisition (reuse "Download from website")
$ws.Range("F4").Value = "Download from website"

# G4: Initial Observations/Challenges
$ws.Range("G4").Value = "Scheme-wise data for Annadata Sukhibhava (current farmer scheme) and Vaddi Leni Runalu (VLR) found in 'LIST OF SCHEMES' table starting on Page 67. All figures are in Rupees Lakhs. YSR Rythu Bharosa was the previous scheme, Annadata Sukhibhava is the current one"
$ws.Range("G4").WrapText = $true

# H4: Local File Path
$ws.Range("H4").Value = """C:\Project_AP_Welfare_Dashboard\Data\Raw\Agriculture_data.pdf"""
$ws.Range("H4").WrapText = $true

# Column G got wider to fit the new observations text; row 4 grows to fit the wrapped text
$ws.Columns.Item(7).ColumnWidth = 34
$ws.Rows.Item(4).RowHeight = 100.8

# Match final selection/active cell as left by the editor
$ws.Range("H4").Select() | Out-Null

# Page orientation touched during this edit
$ws.PageSetup.Orientation = 1

$wb.Save()
